$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44757

# Row 3
$ws.Range("D3").Value = 44812
$ws.Range("J3").Value = 80

# Row 4
$ws.Range("D4").Value = 44839
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 16000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 16000
$ws.Range("P4").Value = 1067

# Row 5
$ws.Range("D5").Value = 44819
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 20000
$ws.Range("L5").Value = 20000
$ws.Range("M5").Value = 20000
$ws.Range("P5").Value = 1333

# Row 7
$ws.Range("D7").Value = 44825
$ws.Range("J7").Value = 30

# Row 8
$ws.Range("D8").Value = 44771
$ws.Range("J8").Value = 40
$ws.Range("K8").Value = 20000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 20000
$ws.Range("P8").Value = 1333

# Row 9
$ws.Range("D9").Value = 44830
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 12000
$ws.Range("P9").Value = 800

# Row 10
$ws.Range("D10").Value = 44749
$ws.Range("K10").Value = 20000
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = 20000
$ws.Range("P10").Value = 1333

# Row 11
$ws.Range("D11").Value = 44508
$ws.Range("J11").Value = 40
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 10000
$ws.Range("P11").Value = 667

# Row 12
$ws.Range("D12").Value = 44827
$ws.Range("J12").Value = 20

# Row 13
$ws.Range("D13").Value = 45134
$ws.Range("J13").Value = 5
$ws.Range("K13").Value = 20000
$ws.Range("L13").Value = 20000
$ws.Range("M13").Value = 20000
$ws.Range("P13").Value = 1333

# Row 15
$ws.Range("D15").Value = 44769
$ws.Range("J15").Value = 50

# Row 16
$ws.Range("D16").Value = 44838
$ws.Range("J16").Value = 10
$ws.Range("K16").Value = 20000
$ws.Range("L16").Value = 20000
$ws.Range("M16").Value = 20000
$ws.Range("P16").Value = 1333

# Row 17
$ws.Range("D17").Value = 44767
$ws.Range("J17").Value = 50

# Row 18
$ws.Range("D18").Value = 44837
$ws.Range("J18").Value = 80

# Row 19
$ws.Range("D19").Value = 44841
$ws.Range("J19").Value = 20
$ws.Range("K19").Value = 16000
$ws.Range("L19").Value = 16000
$ws.Range("M19").Value = 16000
$ws.Range("P19").Value = 1067

# Row 20
$ws.Range("D20").Value = 44811

# Row 21
$ws.Range("D21").Value = 44813
$ws.Range("J21").Value = 20

# Row 22
$ws.Range("D22").Value = 44776
$ws.Range("J22").Value = 80

# Row 23
$ws.Range("D23").Value = 44518
$ws.Range("K23").Value = 10000
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = 10000
$ws.Range("P23").Value = 667

# Row 24
$ws.Range("D24").Value = 44756
$ws.Range("J24").Value = 80

# Row 25
$ws.Range("D25").Value = 44826
$ws.Range("J25").Value = 50

# Row 26
$ws.Range("D26").Value = 44845
$ws.Range("K26").Value = 16000
$ws.Range("L26").Value = 16000
$ws.Range("M26").Value = 16000
$ws.Range("P26").Value = 1067

# Row 27
$ws.Range("D27").Value = 44525
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = 8000
$ws.Range("L27").Value = 8000
$ws.Range("M27").Value = 8000
$ws.Range("P27").Value = 533
